$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# Row 2: Mining and quarrying -> All other industry
$ws.Cells.Item(2, 2).Value = "All other industry"
$ws.Cells.Item(2, 3).Value = 220
$ws.Cells.Item(2, 4).Value = 143.7947882736156

# Row 3: Transport; storage and communication -> Warehouses and storage
$ws.Cells.Item(3, 2).Value = "Warehouses and storage"
$ws.Cells.Item(3, 3).Value = 220
$ws.Cells.Item(3, 4).Value = 111.1441368078176

# Row 4: Manufacturing -> Manufacturing and light industry
$ws.Cells.Item(4, 2).Value = "Manufacturing and light industry"
$ws.Cells.Item(4, 3).Value = 220
$ws.Cells.Item(4, 4).Value = 133.5708469055375

# Row 5: Mining and quarrying -> All other industry
$ws.Cells.Item(5, 2).Value = "All other industry"
$ws.Cells.Item(5, 3).Value = 480
$ws.Cells.Item(5, 4).Value = 143.7947882736156

# Row 6: Transport; storage and communication -> Warehouses and storage
$ws.Cells.Item(6, 2).Value = "Warehouses and storage"
$ws.Cells.Item(6, 3).Value = 480
$ws.Cells.Item(6, 4).Value = 111.1441368078176

# Row 7: Manufacturing -> Manufacturing and light industry
$ws.Cells.Item(7, 2).Value = "Manufacturing and light industry"
$ws.Cells.Item(7, 3).Value = 480
$ws.Cells.Item(7, 4).Value = 133.5708469055375
